$d = $word.ActiveDocument

# --- Add new custom paragraph styles (mirrors styles.xml / stylesWithEffects.xml) ---

# tei_signed
$teiSigned = $d.Styles.Add("teisigned", 1)
$teiSigned.BaseStyle = $d.Styles("Normal")
$teiSigned.NameLocal = "tei_signed"
$teiSigned.QuickStyle = $true
$teiSigned.ParagraphFormat.SpaceBefore = 18
$teiSigned.ParagraphFormat.LeftIndent = 21.55
$teiSigned.ParagraphFormat.FirstLineIndent = -21.55

# tei_speech
$teiSpeech = $d.Styles.Add("teispeech", 1)
$teiSpeech.BaseStyle = $d.Styles("Normal")
$teiSpeech.NameLocal = "tei_speech"
$teiSpeech.QuickStyle = $true
$teiSpeech.ParagraphFormat.LeftIndent = 21.6
$teiSpeech.ParagraphFormat.FirstLineIndent = -21.6

# GeneratedTitle
$genTitle = $d.Styles.Add("GeneratedTitle", 1)
$genTitle.BaseStyle = $d.Styles("Title")
$genTitle.NameLocal = "GeneratedTitle"
$genTitle.QuickStyle = $true

# GeneratedSubTitle
$genSubTitle = $d.Styles.Add("GeneratedSubTitle", 1)
$genSubTitle.BaseStyle = $d.Styles("Subtitle")
$genSubTitle.NameLocal = "GeneratedSubTitle"
$genSubTitle.QuickStyle = $true

# --- Apply the GeneratedSubTitle style to the first paragraph, replacing its
#     direct indent formatting ---
$p = $d.Paragraphs(1)
$p.Style = $d.Styles("GeneratedSubTitle")
